$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing totals (July 1 and July 15)
$ws.Range("B2").Value = 18036.68
$ws.Range("B12").Value = 47778.6

# Insert a new row before row 13, shifting existing rows (13-74) down to (14-75)
$ws.Rows.Item(13).Insert($xlShiftDown)

# Fill in the newly inserted row 13 with the new data point (Day 16, July 2025)
$ws.Range("A13").Value = 16
$ws.Range("B13").Value = 15570.15
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 2025
$ws.Range("E13").Value = "07/2025"
